$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOSPITALES")

# --- Update existing row 97 (Puerto Lempira hospital) ---
$ws.Cells.Item(97, 21).Value = "Hospital Puerto Lempira"
$ws.Cells.Item(97, 22).Value = 15.265726000000001
$ws.Cells.Item(97, 23).Value = -83.781515999999996

# --- Append new rows 114-121 ---
# columns: G=7 Admin1nombre, K=11 Admin2nombre, S=19 Subcategoria, U=21 Nombre, V=22 Latitud, W=23 Longitud

# Row 114
$ws.Cells.Item(114, 7).Value = "Gracias a Dios"
$ws.Cells.Item(114, 19).Value = "centro de salud"
$ws.Cells.Item(114, 21).Value = "Cesamo de Ibans"
$ws.Cells.Item(114, 22).Value = 15.908462
$ws.Cells.Item(114, 23).Value = -84.814170000000004

# Row 115
$ws.Cells.Item(115, 7).Value = "Santa Barbará"
$ws.Cells.Item(115, 11).Value = "Santa Barbará"
$ws.Cells.Item(115, 19).Value = "hospital"
$ws.Cells.Item(115, 21).Value = "Hospital Santa Barbará Integrado"
$ws.Cells.Item(115, 22).Value = 14.925007000000001
$ws.Cells.Item(115, 23).Value = -88.237927999999997

# Row 116
$ws.Cells.Item(116, 7).Value = "Santa Barbará"
$ws.Cells.Item(116, 11).Value = "Trinidad"
$ws.Cells.Item(116, 19).Value = "clínica"
$ws.Cells.Item(116, 21).Value = "Clínica de Emergencia Trinidad"
$ws.Cells.Item(116, 22).Value = 15.143121000000001
$ws.Cells.Item(116, 23).Value = -88.237575000000007

# Row 117
$ws.Cells.Item(117, 7).Value = "Santa Barbará"
$ws.Cells.Item(117, 11).Value = "Sula"
$ws.Cells.Item(117, 19).Value = "hospital"
$ws.Cells.Item(117, 21).Value = "Hospital Sula Socorro de lo Atlto"
$ws.Cells.Item(117, 22).Value = 15.247275
$ws.Cells.Item(117, 23).Value = -88.552238000000003

# Row 118
$ws.Cells.Item(118, 7).Value = "Santa Barbará"
$ws.Cells.Item(118, 11).Value = "Petoa"
$ws.Cells.Item(118, 19).Value = "hospital"
$ws.Cells.Item(118, 21).Value = "Hospital Luz de Vida"
$ws.Cells.Item(118, 22).Value = 15.273072000000001
$ws.Cells.Item(118, 23).Value = -88.284302999999994

# Row 119
$ws.Cells.Item(119, 7).Value = "Santa Barbará"
$ws.Cells.Item(119, 11).Value = "Quimistán"
$ws.Cells.Item(119, 19).Value = "hospital"
$ws.Cells.Item(119, 21).Value = "Centro Médico Integral de Occidente"
$ws.Cells.Item(119, 22).Value = 15.348732999999999
$ws.Cells.Item(119, 23).Value = -88.403704000000005

# Row 120
$ws.Cells.Item(120, 7).Value = "El Paraíso"
$ws.Cells.Item(120, 11).Value = "El Paraíso"
$ws.Cells.Item(120, 19).Value = "hospital"
$ws.Cells.Item(120, 21).Value = "Centro Médico San Francisco"
$ws.Cells.Item(120, 22).Value = 13.861846
$ws.Cells.Item(120, 23).Value = -86.554258000000004

# Row 121
$ws.Cells.Item(121, 7).Value = "El Paraíso"
$ws.Cells.Item(121, 11).Value = "El Paraíso"
$ws.Cells.Item(121, 19).Value = "hospital"
$ws.Cells.Item(121, 21).Value = "Hospital Alivio del Sufrimiento"
$ws.Cells.Item(121, 22).Value = 13.865371
$ws.Cells.Item(121, 23).Value = -86.562415000000001

# --- Match final cursor/selection position (W123, two rows below the last data row) ---
$ws.Range("W123").Select()
